$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mississippi's row (row 42) failed on this run: clear the scraped
# numeric/date values (including the date cell's number format), flip
# the two boolean flags to FALSE, and record the connection-reset error
# message in the status column.
$ws.Range("B42:H42").ClearContents()
$ws.Range("B42").ClearFormats()
$ws.Range("I42").Value = $false
$ws.Range("J42").Value = $false
$ws.Range("O42").Value = "An error occurred. ... ConnectionError(ProtocolError('Connection aborted.', ConnectionResetError(104, 'Connection reset by peer')))"
